$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Lakásár Változás"

# Update the changed values in column B
$ws.Range("B5").Value = 0.5
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 6.5
$ws.Range("B8").Value = 8
$ws.Range("B10").Value = 12.5
$ws.Range("B11").Value = 14
$ws.Range("B12").Value = 9
$ws.Range("B13").Value = 11
$ws.Range("B14").Value = 15
$ws.Range("B15").Value = 4
